# Adds one new weekly price record for "Choclo" (Macroferia Regional de Talca)
# by inserting a new row at row 160 - this shifts the existing rows 160:171
# down to 161:172 (dimension grows from A1:R171 to A1:R172) - and fills the
# new row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 160, pushing rows 160-171 down to 161-172.
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row 160 with the new record.
$ws.Cells.Item(160, 1).Value2 = 5
$ws.Cells.Item(160, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(160, 3).Value2 = "Maule"
$ws.Cells.Item(160, 4).Value2 = 44578
$ws.Cells.Item(160, 5).Value2 = 7
$ws.Cells.Item(160, 6).Value2 = 100112024
$ws.Cells.Item(160, 7).Value2 = "Choclo"
$ws.Cells.Item(160, 8).Value2 = "Choclero"
$ws.Cells.Item(160, 9).Value2 = "Primera"
$ws.Cells.Item(160, 10).Value2 = 40000
$ws.Cells.Item(160, 11).Value2 = 200
$ws.Cells.Item(160, 12).Value2 = 200
$ws.Cells.Item(160, 13).Value2 = 200
$ws.Cells.Item(160, 14).Value2 = "`$/unidad"
$ws.Cells.Item(160, 15).Value2 = "Región del Maule"
$ws.Cells.Item(160, 16).Value2 = 200
$ws.Cells.Item(160, 17).Value2 = 1
$ws.Cells.Item(160, 18).Value2 = "Hortaliza"
